$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Samples" tab query (cell B3) is being trimmed down to the "All studies"
# variant requested in the commit: the Tumor / Analyte Type columns are removed
# from the SELECT list, everything else about the query stays the same.
$newSampleQuery = @"
SELECT
    DISTINCT (smp.sample_id) AS "Sample ID",
    sp.participant_id AS "Participant ID", 
    s.study_name AS "Study Name",
    s.phs_accession AS Accession
FROM 
    df_participant sp
JOIN 
    df_study s ON sp."study.phs_accession" = s.phs_accession
JOIN 
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
JOIN
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
WHERE 
  s.phs_accession = 'phs001819' AND gi.library_selection = 'Random'
ORDER BY 
    smp.sample_id ASC
LIMIT 100;
"@

$ws.Range("B3").Value = $newSampleQuery

# Re-point the view at the new query cell (matches the saved selection /
# scroll position recorded in the workbook after the edit).
$ws.Range("B3").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
